# Updated cryptos list (Price / Volume(1h) columns) per upstream commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold plain-text numbers (e.g. "24.124.26", "1.005").
# Excel auto-converts numeric-looking text typed into .Value, so for those
# we briefly force a text format, assign, then clear the format again -
# this mirrors typing into a "@"-formatted cell and matches the source
# workbook, where these cells are stored as inline/shared strings with no
# explicit cell style.

$ws.Range("D2").Value = '24.124.26'
$ws.Range("E2").Value = '  -3.37%  '
$ws.Range("D3").Value = '1.644.69'
$ws.Range("E3").Value = '  -3.28%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.95'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.44%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3906'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.80%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3859'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -4.10%  '
$ws.Range("E9").Value = '  +0.36%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.359'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -7.14%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '48.88'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -7.67%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08462'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -3.80%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '24.10'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -7.46%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.156'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -4.14%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001285'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -4.79%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.504'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -5.57%  '
$ws.Range("D17").Value = '1.645.76'
$ws.Range("E17").Value = '  -3.72%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '94.43'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.70%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06943'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -3.53%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '20.96'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.60%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.954'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -5.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.004'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.21%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.72'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -4.62%  '
$ws.Range("D24").Value = '24.123.13'
$ws.Range("E24").Value = '  -3.37%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.349'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.50%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.727'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -7.24%  '
$ws.Range("E27").Value = '  -5.31%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.992'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +7.65%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '157.81'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.62%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '141.63'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -6.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.395'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -13.04%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.469'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -6.05%  '
$ws.Range("D33").Value = '1.826.98'
$ws.Range("E33").Value = '  -3.73%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.220'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.97%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08029'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -6.38%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9828'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -5.36%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02938'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -6.36%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2711'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -5.81%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.09260'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -3.31%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.478'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.43%  '
$ws.Range("E41").Value = '  -7.96%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7634'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -7.43%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '13.12'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -6.44%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.02'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -6.71%  '
$ws.Range("E45").Value = '  -7.15%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6894'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -6.65%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.093'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -3.53%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.003'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.18%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08424'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -4.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '134.13'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -3.67%  '
$ws.Range("E51").Value = '  -9.43%  '
